$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '67.131.54'
$ws.Range("E2").Value = '  +0.10%  '
$ws.Range("D3").Value = '2.468.43'
$ws.Range("E4").Value = '  +0.08%  '
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '582.49'
$c.Style = "Normal"
$ws.Range("E5").Value = '  -0.17%  '
$ws.Range("E6").Value = '  +3.50%  '
$ws.Range("E7").Value = '  -0.03%  '
$ws.Range("E8").Value = '  -0.41%  '
$ws.Range("E9").Value = '  +2.31%  '
$ws.Range("E10").Value = '  +0.27%  '
$ws.Range("E11").Value = '  +0.19%  '
$ws.Range("E12").Value = '  +1.20%  '
$ws.Range("D13").Value = '2.915.88'
$ws.Range("E13").Value = '  -0.04%  '
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = '25.39'
$c.Style = "Normal"
$ws.Range("E14").Value = '  -0.65%  '
$ws.Range("D15").Value = '66.983.11'
$ws.Range("E15").Value = '  +0.52%  '
$ws.Range("E16").Value = '  +0.29%  '
$ws.Range("D17").Value = '2.460.66'
$ws.Range("E17").Value = '  +0.33%  '
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = '10.92'
$c.Style = "Normal"
$ws.Range("E18").Value = '  -1.45%  '
$ws.Range("E19").Value = '  -1.51%  '
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = '348.00'
$c.Style = "Normal"
$ws.Range("E20").Value = '  -1.59%  '
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = '4.02'
$c.Style = "Normal"
$ws.Range("E21").Value = '  -0.18%  '
$ws.Range("E22").Value = '  +0.11%  '
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = '69.36'
$c.Style = "Normal"
$ws.Range("E23").Value = '  +0.63%  '
$ws.Range("E24").Value = '  -1.31%  '
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = '1.79'
$c.Style = "Normal"
$ws.Range("E25").Value = '  +0.08%  '
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = '9.22'
$c.Style = "Normal"
$ws.Range("E26").Value = '  +0.00%  '
$ws.Range("D27").Value = '2.595.40'
$ws.Range("E27").Value = '  +0.33%  '
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = '0.999'
$c.Style = "Normal"
$ws.Range("E28").Value = '  +0.09%  '
$ws.Range("D29").Value = '0.0₃0900'
$ws.Range("E29").Value = '  -0.21%  '
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = '499.24'
$c.Style = "Normal"
$ws.Range("E30").Value = '  -2.94%  '
$ws.Range("E31").Value = '  -0.15%  '
$ws.Range("E32").Value = '  -0.27%  '
$ws.Range("E33").Value = '  -0.91%  '
$ws.Range("E34").Value = '  -0.01%  '
$ws.Range("E35").Value = '  +2.41%  '
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = '161.46'
$c.Style = "Normal"
$ws.Range("E36").Value = '  +1.77%  '
$ws.Range("E37").Value = '  +0.10%  '
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = '18.16'
$c.Style = "Normal"
$ws.Range("E38").Value = '  -0.88%  '
$ws.Range("E39").Value = '  -1.37%  '
$ws.Range("E40").Value = '  +0.03%  '
$ws.Range("E41").Value = '  +0.92%  '
$ws.Range("E42").Value = '  +0.23%  '
$ws.Range("E43").Value = '  +0.12%  '
$ws.Range("E44").Value = '  +0.64%  '
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = '142.61'
$c.Style = "Normal"
$ws.Range("E45").Value = '  +1.58%  '
$ws.Range("E46").Value = '  +0.84%  '
$ws.Range("D47").Value = '0.0₆0256'
$ws.Range("E47").Value = '  +0.56%  '
$ws.Range("E48").Value = '  -0.80%  '
$ws.Range("E49").Value = '  +1.38%  '
$ws.Range("E50").Value = '  -1.20%  '
$ws.Range("E51").Value = '  +0.13%  '
